$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (Through 2022-12-04 -> Through 2022-12-05)
$ws.Name = "Through 2022-12-05"

# Update header label in I1 (2022 (through 12-04) -> 2022 (through 12-05))
$ws.Range("I1").Value = "2022 (through 12-05)"

# Update data for November (row 12), December (row 13) and Total (row 14)
$ws.Range("I12").Value = 117
$ws.Range("I13").Value = 19
$ws.Range("I14").Value = 1534
